# Apply the cryptocurrency price/volume refresh described by the commit diff.
# A leading apostrophe forces Excel to store the value as text (matching the
# original inlineStr cells) instead of auto-converting numeric-looking strings
# (e.g. "260.55", "0.120") into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.444.60"
$ws.Range("E2").Value = "'  +0.76%  "
$ws.Range("D3").Value = "'2.012.75"
$ws.Range("E3").Value = "'  -0.25%  "
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("D5").Value = "'260.55"
$ws.Range("E5").Value = "'  +5.45%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "'  -1.92%  "
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("D8").Value = "'56.66"
$ws.Range("E8").Value = "'  -5.59%  "
$ws.Range("E9").Value = "'  -1.54%  "
$ws.Range("D10").Value = "'0.0774"
$ws.Range("E10").Value = "'  -4.53%  "
$ws.Range("E11").Value = "'  -2.85%  "
$ws.Range("D12").Value = "'14.29"
$ws.Range("E12").Value = "'  -5.61%  "
$ws.Range("D13").Value = "'2.308.76"
$ws.Range("E13").Value = "'  -0.20%  "
$ws.Range("D14").Value = "'21.14"
$ws.Range("E14").Value = "'  -5.64%  "
$ws.Range("D15").Value = "'0.803"
$ws.Range("E15").Value = "'  -5.79%  "
$ws.Range("E16").Value = "'  -4.58%  "
$ws.Range("D17").Value = "'1.998.42"
$ws.Range("E17").Value = "'  -1.03%  "
$ws.Range("D18").Value = "'37.315.35"
$ws.Range("E18").Value = "'  +0.65%  "
$ws.Range("D19").Value = "'70.18"
$ws.Range("E19").Value = "'  -0.46%  "
$ws.Range("D20").Value = "'0.0₃0839"
$ws.Range("E20").Value = "'  -3.25%  "
$ws.Range("B21").Value = "'BitcoinCash"
$ws.Range("C21").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'235.85"
$ws.Range("E21").Value = "'  +2.33%  "
$ws.Range("B22").Value = "'Uniswap"
$ws.Range("C22").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.15"
$ws.Range("E22").Value = "'  -1.67%  "
$ws.Range("D23").Value = "'2.63"
$ws.Range("E23").Value = "'  +4.85%  "
$ws.Range("E24").Value = "'  -0.08%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "'  -0.66%  "
$ws.Range("D26").Value = "'165.00"
$ws.Range("E26").Value = "'  +0.87%  "
$ws.Range("D27").Value = "'8.91"
$ws.Range("E27").Value = "'  -5.50%  "
$ws.Range("D28").Value = "'19.71"
$ws.Range("E28").Value = "'  -0.25%  "
$ws.Range("E29").Value = "'  -4.45%  "
$ws.Range("D30").Value = "'1.35"
$ws.Range("E30").Value = "'  -1.71%  "
$ws.Range("D31").Value = "'0.120"
$ws.Range("D32").Value = "'4.61"
$ws.Range("E32").Value = "'  -4.41%  "
$ws.Range("D33").Value = "'0.0645"
$ws.Range("E33").Value = "'  -2.78%  "
$ws.Range("E34").Value = "'  +0.66%  "
$ws.Range("E35").Value = "'  -3.56%  "
$ws.Range("E36").Value = "'  +0.57%  "
$ws.Range("E37").Value = "'  +0.06%  "
$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "'  -3.70%  "
$ws.Range("D39").Value = "'5.36"
$ws.Range("E39").Value = "'  -1.43%  "
$ws.Range("D40").Value = "'3.05"
$ws.Range("E40").Value = "'  +3.94%  "
$ws.Range("D41").Value = "'1.20"
$ws.Range("E41").Value = "'  +0.52%  "
$ws.Range("E42").Value = "'  -0.86%  "
$ws.Range("D43").Value = "'0.0928"
$ws.Range("E43").Value = "'  -5.92%  "
$ws.Range("D44").Value = "'1.415.48"
$ws.Range("E44").Value = "'  +1.89%  "
$ws.Range("D45").Value = "'15.83"
$ws.Range("E45").Value = "'  -5.68%  "
$ws.Range("D46").Value = "'89.88"
$ws.Range("E46").Value = "'  -2.33%  "
$ws.Range("D47").Value = "'1.03"
$ws.Range("E47").Value = "'  -2.85%  "
$ws.Range("B48").Value = "'FraxShare"
$ws.Range("C48").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.06"
$ws.Range("E48").Value = "'  -5.60%  "
$ws.Range("B49").Value = "'MXToken"
$ws.Range("C49").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.92"
$ws.Range("E49").Value = "'  +2.48%  "
$ws.Range("D50").Value = "'2.200.06"
$ws.Range("E50").Value = "'  -0.21%  "
$ws.Range("E51").Value = "'  -9.88%  "
